$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 108 (pushing existing rows down by 2)
$ws.Rows.Item(108).Resize(2).Insert()

# Row 108: new weekly "Primera" quality entry for Chirimoya, Provincia del Elquí
$ws.Range("A108").Value = 3
$ws.Range("B108").Value = "Femacal de La Calera"
$ws.Range("C108").Value = "Coquimbo"
$ws.Range("D108").Value = 44824
$ws.Range("E108").Value = 5
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100107
$ws.Range("H108").Value = "Otros"
$ws.Range("I108").Value = 100107002
$ws.Range("J108").Value = "Chirimoya"
$ws.Range("K108").Value = "Cultivar IV Región"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 30
$ws.Range("N108").Value = 28000
$ws.Range("O108").Value = 28000
$ws.Range("P108").Value = 28000
$ws.Range("Q108").Value = "`$/bandeja 10 kilos"
$ws.Range("R108").Value = "Provincia del Elquí"
$ws.Range("S108").Value = 2800
$ws.Range("T108").Value = 10

# Row 109: new weekly "Segunda" quality entry for Chirimoya, Provincia del Elquí
$ws.Range("A109").Value = 3
$ws.Range("B109").Value = "Femacal de La Calera"
$ws.Range("C109").Value = "Coquimbo"
$ws.Range("D109").Value = 44824
$ws.Range("E109").Value = 5
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100107
$ws.Range("H109").Value = "Otros"
$ws.Range("I109").Value = 100107002
$ws.Range("J109").Value = "Chirimoya"
$ws.Range("K109").Value = "Cultivar IV Región"
$ws.Range("L109").Value = "Segunda"
$ws.Range("M109").Value = 35
$ws.Range("N109").Value = 25000
$ws.Range("O109").Value = 25000
$ws.Range("P109").Value = 25000
$ws.Range("Q109").Value = "`$/bandeja 10 kilos"
$ws.Range("R109").Value = "Provincia del Elquí"
$ws.Range("S109").Value = 2500
$ws.Range("T109").Value = 10
